$wb = $excel.ActiveWorkbook

# First worksheet holds the ValueSet metadata as Property/Value rows.
$ws = $wb.Worksheets.Item(1)

# Row 7 = "Experimental" property. FHIR ValueSets require this boolean
# element; the sheet stores it as literal text "true" (not an Excel
# TRUE boolean), so write it through a text formula and freeze the
# result back to a plain value to dodge Excel's auto boolean coercion
# that a direct .Value assignment of "true" would trigger.
$ws.Range("B7").Formula = "=""true"""
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)

# Row 8 = "Date" property; bump the timestamp to reflect the re-generation.
$ws.Range("B8").Value = "2023-02-01T09:05:11-06:00"
